# Natmi following Dr Hou advice
# Expand the LR-pairs table (Vtn -> Tnfrsf11b) from 2 data rows to 6 data rows,
# adding the "ECs" sending-cluster cluster alongside the existing FAPs / sCs
# clusters, and refreshing all of the computed NATMI statistics.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The table grows from 2 data rows (row 2-3) to 6 data rows (row 2-7), so
# insert 4 fresh rows right after the existing row 2 to make room before
# writing the full, refreshed data set into rows 2-7.
$ws.Rows.Item(3).Insert()
$ws.Rows.Item(3).Insert()
$ws.Rows.Item(3).Insert()
$ws.Rows.Item(3).Insert()

# Columns: A Sending cluster | B Ligand symbol | C Receptor symbol | D Target cluster
# E Ligand-expressing cells | F Ligand detection rate | G Ligand average expression value
# H Ligand total expression value | I Ligand derived specificity of average expression value
# J Ligand derived specificity of total expression value | K Receptor-expressing cells
# L Receptor detection rate | M Receptor average expression value | N Receptor total expression value
# O Receptor derived specificity of average expression value | P Receptor derived specificity of total expression value
# Q Edge average expression weight | R Edge total expression weight
# S Edge average expression derived specificity | T Edge total expression derived specificity
$rows = @(
    @("ECs",  "Vtn", "Tnfrsf11b", "FAPs", 2, 0.6666666666666666, 6.597131,          19.791393,          0.1209543635982448, 0.1209543635982448, 3, 1,                  3.776574666666666,   11.329724, 0.9855052394405499, 0.9855052394405499, 24.91455780728133, 224.231020265532, 0.1192011590592676,  0.1192011590592676),
    @("ECs",  "Vtn", "Tnfrsf11b", "sCs",  2, 0.6666666666666666, 6.597131,          19.791393,          0.1209543635982448, 0.1209543635982448, 1, 0.3333333333333333, 0.05554566666666667, 0.166637,  0.01449476055945007, 0.01449476055945008, 0.3664420394823333, 3.297978355341,    0.001753204538977222, 0.001753204538977222),
    @("FAPs", "Vtn", "Tnfrsf11b", "FAPs", 3, 1,                  21.05317333333333, 63.15952,           0.3859970617919927, 0.3859970617919927, 3, 1,                  3.776574666666666,   11.329724, 0.9855052394405499, 0.9855052394405499, 79.50888106360888, 715.57992957248,  0.3804021268046665,  0.3804021268046665),
    @("FAPs", "Vtn", "Tnfrsf11b", "sCs",  3, 1,                  21.05317333333333, 63.15952,           0.3859970617919927, 0.3859970617919927, 1, 0.3333333333333333, 0.05554566666666667, 0.166637,  0.01449476055945007, 0.01449476055945008, 1.169412548248889,  10.52471293424,   0.005594934987326189, 0.00559493498732619),
    @("sCs",  "Vtn", "Tnfrsf11b", "FAPs", 3, 1,                  26.89201066666667, 80.67603200000001,  0.4930485746097625, 0.4930485746097625, 3, 1,                  3.776574666666666,   11.329724, 0.9855052394405499, 0.9855052394405499, 101.5596862194631, 914.037175975168, 0.4859019535766158,  0.4859019535766158),
    @("sCs",  "Vtn", "Tnfrsf11b", "sCs",  3, 1,                  26.89201066666667, 80.67603200000001,  0.4930485746097625, 0.4930485746097625, 1, 0.3333333333333333, 0.05554566666666667, 0.166637,  0.01449476055945007, 0.01449476055945008, 1.493734660487111,  13.443611944384,  0.007146621033146662, 0.007146621033146663)
)

$startRow = 2
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $row = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    for ($c = 5; $c -le 20; $c++) {
        $ws.Cells.Item($r, $c).Value = $row[$c - 1]
    }
}
